$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the beneficiary-delete query in B2: parameterize the account number
# (was hardcoded '06047900194203', now '{account_number}')
$ws.Range("B2").Value = "BEGIN UPDATE  DC_FUND_TRANSFER_BENEFICIARY TF SET TF.IS_DELETED = 1 WHERE TF.ACCOUNT_NO = '{account_number}' AND TF.CUSTOMER_INFO_ID = (SELECT CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO DCI WHERE DCI.CUSTOMER_NAME = '{customer_name}');DELETE FROM DC_DATA_CACHE DC WHERE DC.CUSTOMER_INFO_ID = (SELECT CI.CUSTOMER_INFO_ID FROM DC_CUSTOMER_INFO CI WHERE CI.CUSTOMER_NAME = '{customer_name}');COMMIT;END;"

# Update the view so that B1 is the top-left visible cell and B2 is selected
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select()
